$d = $word.ActiveDocument

# --- 1) Delete the last paragraph entirely ---
# ("Replace our logic for copying individual files with robocopy...")
$lastParaIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastParaIndex)
$lastPara.Range.Delete()

# --- 2) Move the "_GoBack" bookmark from the start of the
#        "Need to have it so if a drive..." paragraph to its end
#        (right after the last run, before the paragraph mark). ---
$targetParaIndex = $d.Paragraphs.Count
$targetPara = $d.Paragraphs.Item($targetParaIndex)

# Remove the existing (misplaced) bookmark if present.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Position right before the paragraph mark of the target paragraph.
$endPos = $targetPara.Range.End - 1

# Workaround: adding a bookmark directly on a zero-length range located
# immediately before a paragraph mark mis-resolves its position, so we
# temporarily insert a placeholder character, bookmark the (now
# non-empty) range, and then remove the placeholder again -- leaving a
# correctly positioned, collapsed bookmark.
$tempRange = $d.Range($endPos, $endPos)
$tempRange.InsertAfter("X")
$newBookmark = $d.Bookmarks.Add("_GoBack", $tempRange)
$placeholderRange = $d.Range($endPos, $endPos + 1)
$placeholderRange.Text = ""
